# Adjust Excel download/upload to include folders and parents
#
# Before:  one sheet "Seiten" with columns Titel / Beschreibung
# After:   two sheets -
#            "Ordner" (new, first tab)       - ID / In Ordner / Titel
#            "Seiten" (existing, second tab) - ID / In Ordner / Titel / Beschreibung
#
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Ordner" sheet in front of "Seiten" FIRST - adding a sheet
#    invalidates any worksheet handles fetched beforehand, so every other
#    worksheet reference used below is (re-)fetched after this call.
# ---------------------------------------------------------------------------
$ordner = $wb.Worksheets.Add()
$ordner.Name = "Ordner"

$ordner.Range("A1").Value = "ID"
$ordner.Range("B1").Value = "In Ordner"
$ordner.Range("C1").Value = "Titel"
$ordner.Range("A1:C1").Font.Bold = $true

$ordner.Range("C2").Value = "Titel eines Seitenordners"
$ordner.Range("C3").Value = "Ein weiterer Ordner"

$ordner.Columns.Item(1).ColumnWidth = 22.166666666666668
$ordner.Columns.Item(2).ColumnWidth = 9.498697916666666
$ordner.Columns.Item(3).ColumnWidth = 18.721354166666668

$ordner.Range("C8").Select()

# ---------------------------------------------------------------------------
# 2. "Seiten": make room for the new "ID" (A) and "In Ordner" (B) columns in
#    front of the existing "Titel" (A) / "Beschreibung" (B) columns.
# ---------------------------------------------------------------------------
$seiten = $wb.Worksheets.Item("Seiten")

$seiten.Range("A1").EntireColumn.Insert()   # Titel       -> B, Beschreibung -> C
$seiten.Range("A1").EntireColumn.Insert()   # ID column now A, rest shifts to C/D

$seiten.Range("A1").Value = "ID"
$seiten.Range("B1").Value = "In Ordner"

# Replace the old two sample rows with the new four rows (the ID column is
# left blank - the real IDs are generated on import/export).
$seiten.Range("A2:D3").ClearContents()

$seiten.Range("B2").Value = 1
$seiten.Range("C2").Value = "Beispieltitel"
$seiten.Range("D2").Value = "<p>Die Beschreibung ist das, was auf der Seite angezeigt wird.</p>"

$seiten.Range("B3").Value = 1
$seiten.Range("C3").Value = "Zweite Seite"
$seiten.Range("D3").Value = "<p>Die Beschreibung ist normalerweise als <strong>HTML</strong> formatiert.</p>"

$seiten.Range("B4").Value = 2
$seiten.Range("C4").Value = "Dritte Seite"
$seiten.Range("D4").Value = "Es ist auch möglich, normalen Text zu verwenden. Wir werden das Beste daraus machen."

$seiten.Range("B5").Value = 2
$seiten.Range("C5").Value = "Titel der letzten Seite"
$seiten.Range("D5").Value = "Dies ist der Inhalt, der auf der letzten Seite angezeigt wird."

$seiten.Columns.Item(1).ColumnWidth = 4.498697916666667
$seiten.Columns.Item(2).ColumnWidth = 8.276041666666666
$seiten.Columns.Item(3).ColumnWidth = 31.385416666666668
# column D keeps its original width (80)

$seiten.Range("C9").Select()

# ---------------------------------------------------------------------------
# 3. "Seiten" stays the active tab - (re-)fetch it fresh, then activate.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Seiten").Activate()
